$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Gun Pickup" row (row 20): Time Taken (hrs) went from 2 to 4 hours
$ws.Range("D20").Value = 4

# Update the status note for the "Power ups" group (merged E17:E20) to reflect
# that the gun pickup code is broken
$ws.Range("E17").Value = "Gun Code Broken"

# Reflect the author's last on-screen selection/scroll position
$ws.Range("D21").Select()
